$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark that currently splits the
#    title run ("会  议  纪  " | bookmark | "要") into two runs.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Re-apply the same text to the title paragraph so Word collapses the
#    (now bookmark-free) adjoining runs back into a single run.
$titlePara = $d.Paragraphs(1).Range
$titlePara.Find.Execute("会  议  纪  要", $true, $false, $false, $false, $false,
                         $true, 1, $false, "会  议  纪  要", 2)

# 3) Bump the document number from SE2022-03 to SE2022-04.
$d.Content.Find.Execute("SE2022-03", $true, $false, $false, $false, $false,
                         $true, 1, $false, "SE2022-04", 2)

# 4) Re-insert the "_GoBack" bookmark right after "SE2022-04", before the
#    paragraph mark. A collapsed range exactly at a paragraph end behaves
#    oddly when handed to Bookmarks.Add, so nudge past it with a throwaway
#    character, plant the bookmark, then remove the throwaway character.
$numRange = $d.Content
$numRange.Find.Execute("SE2022-04", $true, $false, $false, $false, $false,
                        $true, 1, $false, "", 0)
$numRange.Collapse(0)
$numRange.InsertAfter("X")

$afterNumber = $d.Range($numRange.Start, $numRange.Start)
$d.Bookmarks.Add("_GoBack", $afterNumber)

$tempChar = $d.Range($numRange.Start, $numRange.Start + 1)
$tempChar.Delete()
